# Updated cryptos list values (price + 1h volume/change%) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.524.58'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '1.906.60'
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4845'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.57%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4078'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08141'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.56%  '
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.55'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Value = '1.904.98'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.032'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.094'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.48'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06738'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001042'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('D21').Value = '29.540.11'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.583'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.63%  '
$ws.Range('E24').Value = '  -2.10%  '
$ws.Range('D25').Value = '2.115.83'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.279'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +10.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.106'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.040'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09561'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.528'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.394'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.552'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02270'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06120'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.173'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5958'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.926'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.37'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1860'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.418'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.282'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07741'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.46'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5574'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.961'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '115.19'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.055'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.56%  '
